$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 44; this pushes the existing rows 44-59 down to 45-60
$ws.Rows("44").Insert()

# Populate the newly inserted row 44 with the new weekly price record
$ws.Range("A44").Value = 10
$ws.Range("B44").Value = "Vega Modelo de Temuco"
$ws.Range("C44").Value = "La Araucanía"
$ws.Range("D44").Value = 44518
$ws.Range("E44").Value = 9
$ws.Range("F44").Value = 300000000
$ws.Range("G44").Value = "Espárragos"
$ws.Range("H44").Value = "Sin especificar"
$ws.Range("I44").Value = "Primera"
$ws.Range("J44").Value = 200
$ws.Range("K44").Value = 1400
$ws.Range("L44").Value = 1400
$ws.Range("M44").Value = 1400
$ws.Range("N44").Value = "$/kilo"
$ws.Range("O44").Value = "Región del Maule"
$ws.Range("P44").Value = 1400
$ws.Range("Q44").Value = 1
$ws.Range("R44").Value = "Hortaliza"
